# Regenerate save_data column G ("K") values (Strike# -> K), per commit:
# "regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals"
#
# The new K values below are the recalculated strikeout counts (s_vals) that
# replace the previous "Strike#" values in column G, rows 2-69.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @(2,0,1,2,0,3,1,1,1,3,2,1,1,2,0,2,3,2,3,0,1,1,1,0,2,0,0,3,0,0,0,1,2,1,0,0,2,1,0,4,1,0,0,1,1,1,0,1,3,1,1,0,2,1,0,0,0,0,2,1,2,2,1,1,1,3,1,1)

$startRow = 2
for ($i = 0; $i -lt $newValues.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 7).Value = $newValues[$i]
}
